# Update gh-pages to output generated at 456a3b4
# Applies "want-to-go" count refresh across sheets, plus a content re-sync
# for three rows (19, 43, 44) on the "全部类型" (all-types) aggregate sheet.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item(1)   # 展览
$wsShow = $wb.Worksheets.Item(2)   # 演出
$wsLife = $wb.Worksheets.Item(3)   # 本地生活
$wsAll  = $wb.Worksheets.Item(4)   # 全部类型

# ---- 展览 (sheet 1): 想去人数 (F column) updates ----
$wsExpo.Range("F3").Value = 835
$wsExpo.Range("F4").Value = 4272
$wsExpo.Range("F7").Value = 3401
$wsExpo.Range("F8").Value = 953
$wsExpo.Range("F11").Value = 274
$wsExpo.Range("F12").Value = 2298
$wsExpo.Range("F13").Value = 1244
$wsExpo.Range("F16").Value = 492
$wsExpo.Range("F17").Value = 240
$wsExpo.Range("F19").Value = 9488
$wsExpo.Range("F20").Value = 5889
$wsExpo.Range("F21").Value = 376
$wsExpo.Range("F22").Value = 198
$wsExpo.Range("F24").Value = 98
$wsExpo.Range("F33").Value = 4767
$wsExpo.Range("F35").Value = 991
$wsExpo.Range("F36").Value = 121

# ---- 演出 (sheet 2): 想去人数 (F column) updates ----
$wsShow.Range("F4").Value = 9

# ---- 本地生活 (sheet 3): 想去人数 (F column) updates ----
$wsLife.Range("F2").Value = 8613
$wsLife.Range("F4").Value = 1488

# ---- 全部类型 (sheet 4): 想去人数 (F column) updates ----
$wsAll.Range("F2").Value = 8613
$wsAll.Range("F5").Value = 1488
$wsAll.Range("F7").Value = 4272
$wsAll.Range("F10").Value = 3401
$wsAll.Range("F11").Value = 953
$wsAll.Range("F13").Value = 9
$wsAll.Range("F14").Value = 274
$wsAll.Range("F15").Value = 2298
$wsAll.Range("F20").Value = 1244
$wsAll.Range("F24").Value = 492
$wsAll.Range("F25").Value = 240
$wsAll.Range("F27").Value = 9488
$wsAll.Range("F30").Value = 376
$wsAll.Range("F31").Value = 198
$wsAll.Range("F33").Value = 98
$wsAll.Range("F42").Value = 4767

# ---- 全部类型 row 19: re-sync to the updated 七夕/Beyond tribute concert listing ----
$wsAll.Range("C19").Value = '杭州·【七夕巨献·早鸟6折】真的爱你”致敬Beyond·黄家驹31周年演唱会·630乐团再现91殿堂级演出'
$wsAll.Range("D19").Value = '湖墅南路136-138号 浙话艺术剧院'
$wsAll.Range("E19").Value = '2024.08.10 19:30-08.10 21:30'
$wsAll.Range("F19").Value = 6
$wsAll.Range("G19").Value = 60
$wsAll.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=85333'
$wsAll.Range("I19").Value = '//i1.hdslb.com/bfs/openplatform/202405/uYt32zt21715221330023.jpeg'

# ---- 全部类型 row 43: was the expired 布谷布 listing, now holds the Eternal only-展 listing ----
# (leading apostrophe keeps the date-looking string as literal text, matching the
#  original inlineStr cell instead of letting Excel coerce it to a date serial)
$wsAll.Range("B43").Value = "'2024-09-16"
$wsAll.Range("C43").Value = '杭州·Eternal时光国乙only展（日+夜场）'
$wsAll.Range("D43").Value = '创意路1号 中国智谷富春园区'
$wsAll.Range("E43").Value = '2024.09.16 09:30-09.17 17:00'
$wsAll.Range("F43").Value = 991
$wsAll.Range("G43").Value = 75
$wsAll.Range("H43").Value = 'https://show.bilibili.com/platform/detail.html?id=89250'
$wsAll.Range("I43").Value = '//i0.hdslb.com/bfs/openplatform/202407/VVNYXGdJ1720966510693.png'

# ---- 全部类型 row 44: new listing for 第五人格only ----
$wsAll.Range("C44").Value = '杭州·第五人格only'
$wsAll.Range("D44").Value = '望江东路333号 杭州瑞莱克斯大酒店'
$wsAll.Range("E44").Value = '2024.09.16 10:00-09.16 17:00'
$wsAll.Range("F44").Value = 121
$wsAll.Range("G44").Value = 60
$wsAll.Range("H44").Value = 'https://show.bilibili.com/platform/detail.html?id=89550'
$wsAll.Range("I44").Value = '//i2.hdslb.com/bfs/openplatform/202407/gFZS33XD1721303396870.jpeg'

